# Add 2022-Q3 data: insert a new latest quarter in front of the existing
# per-quarter sheets (pushing every older quarter down one slot), and append
# a brand-new trailing sheet for the quarter that falls off the back
# (2020-Q4), whose figures used to live in the old last sheet.

$wb = $excel.ActiveWorkbook

# D2/E2/F2/G2 hold numeric-looking text (e.g. "1.20", trailing zero kept),
# so they're written with a leading apostrophe to keep them text instead of
# letting auto-conversion turn them into floats that would drop the zero.
# H2 is a genuine number (rank), so it's left alone.
function Set-QuarterSheet($ws, $d1Label, $d2, $e2, $f2, $g2, $h2) {
    $ws.Range("D1").Value = $d1Label
    $ws.Range("D2").Value = "'" + $d2
    $ws.Range("E2").Value = "'" + $e2
    $ws.Range("F2").Value = "'" + $f2
    $ws.Range("G2").Value = "'" + $g2
    $ws.Range("H2").Value = $h2
}

# Shift every existing per-quarter sheet to the next-newer quarter's
# identity (name + figures), freeing up the oldest slot for the brand-new
# trailing sheet we append at the end.
$s = $wb.Worksheets.Item("2022-Q2")
$s.Name = "2022-Q3"
Set-QuarterSheet $s "基金规模" "1.12" "90.06" "2.80" "0.0314" 6

$s = $wb.Worksheets.Item("2022-Q1")
$s.Name = "2022-Q2"
Set-QuarterSheet $s "基金规模" "1.20" "88.32" "2.24" "0.0269" 10

$s = $wb.Worksheets.Item("2021-Q4")
$s.Name = "2022-Q1"
Set-QuarterSheet $s "基金规模" "1.35" "88.71" "2.53" "0.0342" 7

$s = $wb.Worksheets.Item("2021-Q3")
$s.Name = "2021-Q4"
Set-QuarterSheet $s "基金规模" "1.22" "90.04" "2.70" "0.0329" 6

$s = $wb.Worksheets.Item("2021-Q2")
$s.Name = "2021-Q3"
Set-QuarterSheet $s "基金金额" "1.22" "88.77" "2.83" "0.0345" 7

$s = $wb.Worksheets.Item("2021-Q1")
$s.Name = "2021-Q2"
Set-QuarterSheet $s "基金金额" "1.34" "90.09" "3.22" "0.0431" 5

$s = $wb.Worksheets.Item("2020-Q4")
$s.Name = "2021-Q1"
Set-QuarterSheet $s "基金金额" "1.35" "87.46" "3.27" "0.0441" 6

# Append the new trailing 2020-Q4 sheet (same layout as its siblings),
# carrying the figures that used to sit in the old last sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q4 = $wb.Worksheets.Add($null, $lastSheet)
$q4.Name = "2020-Q4"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金金额"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'007280"
$q4.Range("C2").Value = "上投摩根日本精选股票（QDII）"
$q4.Range("D2").Value = "'1.11"
$q4.Range("E2").Value = "'94.14"
$q4.Range("F2").Value = "'3.83"
$q4.Range("G2").Value = "'0.0425"
$q4.Range("H2").Value = 6

# Apply the same header style used by the other quarter sheets.
$q4.Range("B1:H1").Style = $wb.Worksheets.Item("2021-Q1").Range("B1:H1").Style

# --- "总计" (summary) sheet: insert a new row for 2022-Q3 and append a new
# row for the 2020-Q4 quarter that the shift pushes past the old range. ---
$zj = $wb.Worksheets.Item("总计")

$zj.Rows.Item(2).Insert()
$zj.Range("A3").Copy($zj.Range("A2"))
$zj.Range("B2:D2").ClearFormats()

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0.03

# Column A is a plain 0-based row index -- independent of the data shift --
# so every row below the insertion point needs it re-stamped.
for ($r = 3; $r -le 8; $r++) {
    $zj.Range("A$r").Value = $r - 2
}

# The row that now represents 2021-Q3 drops to 0.03 (matches its sheet).
$zj.Range("D6").Value = 0.03

# Append the trailing row for 2020-Q4 (previously past the sheet's range).
$zj.Range("A8").Copy($zj.Range("A9"))
$zj.Range("A9").Value = 7
$zj.Range("B9").Value = "2020-Q4"
$zj.Range("C9").Value = 1
$zj.Range("D9").Value = 0.04

Write-Output "done"
